$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Returns sheet: rework rows 2-4 from Customer returns into Expense (supplier)
# returns, tweak rows 5-6 customer details, and append two brand-new customer
# return rows (7-8) so the table has data for both the "Customer" and
# "Expense" return-type tabs.
# ---------------------------------------------------------------------------
$returns = $wb.Worksheets.Item("Returns")

# Row 2: RET-001 -> Expense / Supplier return
$returns.Cells.Item(2,2).Value = "Expense"
$returns.Cells.Item(2,3).Value = ""
$returns.Cells.Item(2,4).Value = "SUP-006"
$returns.Cells.Item(2,6).Value = 5
$returns.Cells.Item(2,7).Value = 45995.90923776707
$returns.Cells.Item(2,8).Value = 411.1
$returns.Cells.Item(2,10).Value = "Defective"

# Row 3: RET-002 -> Expense / Supplier return
$returns.Cells.Item(3,2).Value = "Expense"
$returns.Cells.Item(3,3).Value = ""
$returns.Cells.Item(3,4).Value = "SUP-002"
$returns.Cells.Item(3,5).Value = "Microsoft 365 Business Premium"
$returns.Cells.Item(3,7).Value = 46009.90923776707
$returns.Cells.Item(3,8).Value = 446.24
$returns.Cells.Item(3,9).Value = 0

# Row 4: RET-003 -> Expense / Supplier return
$returns.Cells.Item(4,2).Value = "Expense"
$returns.Cells.Item(4,3).Value = ""
$returns.Cells.Item(4,4).Value = "SUP-007"
$returns.Cells.Item(4,5).Value = "Cisco Catalyst 1000-24T Switch"
$returns.Cells.Item(4,7).Value = 46022.90923776707
$returns.Cells.Item(4,8).Value = 170.61
$returns.Cells.Item(4,10).Value = "Changed mind"

# Row 5: RET-004 stays a Customer return, but swaps to a different customer
# and product
$returns.Cells.Item(5,3).Value = "CUS-007"
$returns.Cells.Item(5,5).Value = "Dell OptiPlex 7010 Desktop"
$returns.Cells.Item(5,6).Value = 1
$returns.Cells.Item(5,7).Value = 45999.90923776707
$returns.Cells.Item(5,8).Value = 207.8
$returns.Cells.Item(5,9).Value = 0
$returns.Cells.Item(5,10).Value = "Changed mind"

# Row 6: RET-005 stays a Customer return, but swaps to a different customer
# and product
$returns.Cells.Item(6,3).Value = "CUS-006"
$returns.Cells.Item(6,5).Value = "Ubiquiti UniFi Access Point"
$returns.Cells.Item(6,7).Value = 46016.90923776707
$returns.Cells.Item(6,8).Value = 248.64
$returns.Cells.Item(6,9).Value = 0
$returns.Cells.Item(6,10).Value = "Changed mind"

# Apply the date number format + matching style to the Return Date cells that
# changed in rows 2-6 (mirrors the existing "yyyy-mm-dd h:mm:ss" style used by
# the rest of column G).
for ($r = 2; $r -le 6; $r++) {
    $returns.Cells.Item($r,7).NumberFormat = "yyyy-mm-dd h:mm:ss"
}

# New row 7: RET-006, Customer return
$returns.Cells.Item(7,1).Value = "RET-006"
$returns.Cells.Item(7,2).Value = "Customer"
$returns.Cells.Item(7,3).Value = "CUS-007"
$returns.Cells.Item(7,5).Value = "Monthly IT Support Contract"
$returns.Cells.Item(7,6).Value = 2
$returns.Cells.Item(7,7).Value = 46048.90923776707
$returns.Cells.Item(7,7).NumberFormat = "yyyy-mm-dd h:mm:ss"
$returns.Cells.Item(7,8).Value = 272.62
$returns.Cells.Item(7,9).Value = 0
$returns.Cells.Item(7,10).Value = "Changed mind"
$returns.Cells.Item(7,11).Value = "Completed"

# New row 8: RET-007, Customer return
$returns.Cells.Item(8,1).Value = "RET-007"
$returns.Cells.Item(8,2).Value = "Customer"
$returns.Cells.Item(8,3).Value = "CUS-008"
$returns.Cells.Item(8,5).Value = "Monthly IT Support Contract"
$returns.Cells.Item(8,6).Value = 1
$returns.Cells.Item(8,7).Value = 46028.90923776707
$returns.Cells.Item(8,7).NumberFormat = "yyyy-mm-dd h:mm:ss"
$returns.Cells.Item(8,8).Value = 98.38
$returns.Cells.Item(8,9).Value = 0
$returns.Cells.Item(8,10).Value = "Not as described"
$returns.Cells.Item(8,11).Value = "Completed"

# ---------------------------------------------------------------------------
# Lost Damaged sheet: refresh the "Date Discovered" timestamps (regenerated
# sample data) and drop the stray empty Notes cells in rows 2-5.
# ---------------------------------------------------------------------------
$lostDamaged = $wb.Worksheets.Item("Lost Damaged")

$lostDamaged.Cells.Item(2,5).Value = 46005.89542372686
$lostDamaged.Cells.Item(3,5).Value = 46047.89542372686
$lostDamaged.Cells.Item(4,5).Value = 46045.89542372686
$lostDamaged.Cells.Item(5,5).Value = 46037.89542372686

$lostDamaged.Cells.Item(2,8).ClearContents()
$lostDamaged.Cells.Item(3,8).ClearContents()
$lostDamaged.Cells.Item(4,8).ClearContents()
$lostDamaged.Cells.Item(5,8).ClearContents()
